# Generate Report for Handback
#
# The handback CI run produced a new pair of source/target file UUIDs and a
# new content hash for the .xlf round-trip files, plus refreshed handoff /
# handback timestamps. Update the three status sheets (Overview, zh-cn,
# de-de) so every cell that echoes the old UUIDs / hash / timestamps shows
# the new ones.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "155daaac-08c1-48d0-964f-fe84de3a97db"
$newUuid1 = "3819ca96-5dc3-45a3-bc6a-3ba87bcd3580"

$oldUuid2 = "21a6145f-f2ad-4ad9-ae1b-10e2d89eec22"
$newUuid2 = "ffff2104bfa4-17a4-4de5-833c-a10edac9086d"

$newHash = "f7e8e85763eb69c36778c29e41793ccfa9dddb62"

$newFileMd1 = "$newUuid1.md"
$newFileMd2 = "$newUuid2.md"

$newXlfZhCn = "$newUuid1.$newHash.zh-cn.xlf"
$newXlfDeDe = "$newUuid1.$newHash.de-de.xlf"

$zhHandoffTime = "2016-03-19 16:51:53"
$zhHandbackTime = "2016-03-19 16:52:34"

$deHandoffTime = "2016-03-19 16:51:56"
$deHandbackTime = "2016-03-19 16:52:39"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = $newFileMd1
$overview.Range("A3").Value = $newFileMd2

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = $newFileMd1
$zhcn.Range("D2").Value = $newXlfZhCn
$zhcn.Range("E2").Value = $zhHandoffTime
$zhcn.Range("F2").Value = $newFileMd1
$zhcn.Range("G2").Value = $newXlfZhCn
$zhcn.Range("H2").Value = $zhHandbackTime

$zhcn.Range("A3").Value = $newFileMd2
$zhcn.Range("D3").Value = $newXlfZhCn
$zhcn.Range("E3").Value = $zhHandoffTime
$zhcn.Range("F3").Value = $newFileMd2
$zhcn.Range("G3").Value = $newXlfZhCn
$zhcn.Range("H3").Value = $zhHandbackTime

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = $newFileMd1
$dede.Range("D2").Value = $newXlfDeDe
$dede.Range("E2").Value = $deHandoffTime
$dede.Range("F2").Value = $newFileMd1
$dede.Range("G2").Value = $newXlfDeDe
$dede.Range("H2").Value = $deHandbackTime

$dede.Range("A3").Value = $newFileMd2
$dede.Range("D3").Value = $newXlfDeDe
$dede.Range("E3").Value = $deHandoffTime
$dede.Range("F3").Value = $newFileMd2
$dede.Range("G3").Value = $newXlfDeDe
$dede.Range("H3").Value = $deHandbackTime
